$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '28.090.24'; NumericLike = $true },
    @{ Cell = 'E2'; Value = '  +1.52%  '; NumericLike = $false },
    @{ Cell = 'D3'; Value = '1.895.18'; NumericLike = $true },
    @{ Cell = 'E3'; Value = '  +1.65%  '; NumericLike = $false },
    @{ Cell = 'D4'; Value = '1.017'; NumericLike = $true },
    @{ Cell = 'E4'; Value = '  +1.24%  '; NumericLike = $false },
    @{ Cell = 'D5'; Value = '336.98'; NumericLike = $true },
    @{ Cell = 'E5'; Value = '  +1.77%  '; NumericLike = $false },
    @{ Cell = 'D6'; Value = '1.014'; NumericLike = $true },
    @{ Cell = 'E6'; Value = '  +0.99%  '; NumericLike = $false },
    @{ Cell = 'D7'; Value = '0.4785'; NumericLike = $true },
    @{ Cell = 'E7'; Value = '  +2.24%  '; NumericLike = $false },
    @{ Cell = 'D8'; Value = '0.3978'; NumericLike = $true },
    @{ Cell = 'E8'; Value = '  +1.28%  '; NumericLike = $false },
    @{ Cell = 'D9'; Value = '47.65'; NumericLike = $true },
    @{ Cell = 'E9'; Value = '  -0.05%  '; NumericLike = $false },
    @{ Cell = 'D10'; Value = '0.08050'; NumericLike = $true },
    @{ Cell = 'E10'; Value = '  +0.71%  '; NumericLike = $false },
    @{ Cell = 'D11'; Value = '1.026'; NumericLike = $true },
    @{ Cell = 'E11'; Value = '  +0.39%  '; NumericLike = $false },
    @{ Cell = 'D12'; Value = '22.11'; NumericLike = $true },
    @{ Cell = 'E12'; Value = '  +1.97%  '; NumericLike = $false },
    @{ Cell = 'D13'; Value = '1.894.07'; NumericLike = $true },
    @{ Cell = 'E13'; Value = '  +1.97%  '; NumericLike = $false },
    @{ Cell = 'D14'; Value = '6.042'; NumericLike = $true },
    @{ Cell = 'E14'; Value = '  +2.12%  '; NumericLike = $false },
    @{ Cell = 'D15'; Value = '7.249'; NumericLike = $true },
    @{ Cell = 'E15'; Value = '  +2.01%  '; NumericLike = $false },
    @{ Cell = 'D16'; Value = '1.018'; NumericLike = $true },
    @{ Cell = 'E16'; Value = '  +1.03%  '; NumericLike = $false },
    @{ Cell = 'D17'; Value = '88.92'; NumericLike = $true },
    @{ Cell = 'E17'; Value = '  +2.92%  '; NumericLike = $false },
    @{ Cell = 'D18'; Value = '0.06798'; NumericLike = $true },
    @{ Cell = 'E18'; Value = '  +2.59%  '; NumericLike = $false },
    @{ Cell = 'E19'; Value = '  +0.90%  '; NumericLike = $false },
    @{ Cell = 'D20'; Value = '17.20'; NumericLike = $true },
    @{ Cell = 'E20'; Value = '  -0.01%  '; NumericLike = $false },
    @{ Cell = 'E21'; Value = '  +0.96%  '; NumericLike = $false },
    @{ Cell = 'D22'; Value = '28.092.96'; NumericLike = $true },
    @{ Cell = 'E22'; Value = '  +1.53%  '; NumericLike = $false },
    @{ Cell = 'D23'; Value = '5.553'; NumericLike = $true },
    @{ Cell = 'E23'; Value = '  +1.50%  '; NumericLike = $false },
    @{ Cell = 'D24'; Value = '11.09'; NumericLike = $true },
    @{ Cell = 'D25'; Value = '2.356'; NumericLike = $true },
    @{ Cell = 'E25'; Value = '  +2.02%  '; NumericLike = $false },
    @{ Cell = 'D26'; Value = '2.111.87'; NumericLike = $true },
    @{ Cell = 'E26'; Value = '  +1.61%  '; NumericLike = $false },
    @{ Cell = 'D27'; Value = '160.98'; NumericLike = $true },
    @{ Cell = 'E27'; Value = '  +1.30%  '; NumericLike = $false },
    @{ Cell = 'D28'; Value = '20.16'; NumericLike = $true },
    @{ Cell = 'E28'; Value = '  +0.04%  '; NumericLike = $false },
    @{ Cell = 'D29'; Value = '2.123'; NumericLike = $true },
    @{ Cell = 'E29'; Value = '  +2.29%  '; NumericLike = $false },
    @{ Cell = 'D30'; Value = '5.570'; NumericLike = $true },
    @{ Cell = 'E30'; Value = '  +0.78%  '; NumericLike = $false },
    @{ Cell = 'D31'; Value = '122.44'; NumericLike = $true },
    @{ Cell = 'E31'; Value = '  -0.03%  '; NumericLike = $false },
    @{ Cell = 'D32'; Value = '0.9877'; NumericLike = $true },
    @{ Cell = 'E32'; Value = '  +2.39%  '; NumericLike = $false },
    @{ Cell = 'D33'; Value = '0.09634'; NumericLike = $true },
    @{ Cell = 'E33'; Value = '  +1.65%  '; NumericLike = $false },
    @{ Cell = 'D34'; Value = '3.649'; NumericLike = $true },
    @{ Cell = 'E34'; Value = '  +1.40%  '; NumericLike = $false },
    @{ Cell = 'D35'; Value = '5.385'; NumericLike = $true },
    @{ Cell = 'E35'; Value = '  +1.75%  '; NumericLike = $false },
    @{ Cell = 'D36'; Value = '1.377'; NumericLike = $true },
    @{ Cell = 'E36'; Value = '  -4.64%  '; NumericLike = $false },
    @{ Cell = 'B37'; Value = 'VeChain'; NumericLike = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; NumericLike = $false },
    @{ Cell = 'D37'; Value = '0.02266'; NumericLike = $true },
    @{ Cell = 'E37'; Value = '  +0.91%  '; NumericLike = $false },
    @{ Cell = 'B38'; Value = 'Hedera'; NumericLike = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; NumericLike = $false },
    @{ Cell = 'D38'; Value = '0.06124'; NumericLike = $true },
    @{ Cell = 'E38'; Value = '  +1.15%  '; NumericLike = $false },
    @{ Cell = 'D39'; Value = '1.210'; NumericLike = $true },
    @{ Cell = 'E39'; Value = '  -1.40%  '; NumericLike = $false },
    @{ Cell = 'D40'; Value = '8.241'; NumericLike = $true },
    @{ Cell = 'E40'; Value = '  +1.73%  '; NumericLike = $false },
    @{ Cell = 'E41'; Value = '  +0.97%  '; NumericLike = $false },
    @{ Cell = 'D42'; Value = '0.6007'; NumericLike = $true },
    @{ Cell = 'E42'; Value = '  +0.73%  '; NumericLike = $false },
    @{ Cell = 'D43'; Value = '0.1908'; NumericLike = $true },
    @{ Cell = 'E43'; Value = '  +1.20%  '; NumericLike = $false },
    @{ Cell = 'D44'; Value = '10.42'; NumericLike = $true },
    @{ Cell = 'E44'; Value = '  +2.23%  '; NumericLike = $false },
    @{ Cell = 'D45'; Value = '1.279'; NumericLike = $true },
    @{ Cell = 'E45'; Value = '  +2.24%  '; NumericLike = $false },
    @{ Cell = 'D46'; Value = '0.5702'; NumericLike = $true },
    @{ Cell = 'E46'; Value = '  +0.41%  '; NumericLike = $false },
    @{ Cell = 'D47'; Value = '12.36'; NumericLike = $true },
    @{ Cell = 'E47'; Value = '  +1.99%  '; NumericLike = $false },
    @{ Cell = 'D48'; Value = '1.943'; NumericLike = $true },
    @{ Cell = 'E48'; Value = '  +0.86%  '; NumericLike = $false },
    @{ Cell = 'D49'; Value = '3.388'; NumericLike = $true },
    @{ Cell = 'E49'; Value = '  +0.23%  '; NumericLike = $false },
    @{ Cell = 'D50'; Value = '0.06841'; NumericLike = $true },
    @{ Cell = 'E50'; Value = '  +0.27%  '; NumericLike = $false },
    @{ Cell = 'D51'; Value = '112.93'; NumericLike = $true },
    @{ Cell = 'E51'; Value = '  -0.54%  '; NumericLike = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    if ($u.NumericLike) {
        # Force these to stay text (matches source data, which stores prices
        # as literal strings like "1.017" / "28.090.24"), even though they
        # look numeric. A leading apostrophe stops Excel's autoconvert from
        # turning them into numbers (which would lose trailing/leading
        # zeros), then resetting the style back to Normal drops the
        # quote-prefix formatting so no stray style is left behind.
        $cell.Value = "'" + $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
